$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grade value 5 for the newly completed homework cells
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5

$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 5

$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 5

$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 5

$ws.Range("F24").Value = 5
$ws.Range("F26").Value = 5

# Update the active selection to match the final cursor position
$ws.Range("F20").Select()
